$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 2-147: columns A(date), B, C, D, E, and G label (quarter date+index)
# F (Quarter number) is unchanged from the original file and is not touched here.
$rows = @(
    @(2, 45917, 0.005, 5.693, 0, 0, "17.09.20251"),
    @(3, 45917.01041666666, 0.128, 0.129, 0, 0, "17.09.20252"),
    @(4, 45917.02083333334, 0.123, 0.414, 0, 0, "17.09.20253"),
    @(5, 45917.03125, 0, 7.83, 0, 0, "17.09.20254"),
    @(6, 45917.04166666666, 0, 3.181, 0, 0, "17.09.20255"),
    @(7, 45917.05208333334, 0.002, 3.068, 0, 0, "17.09.20256"),
    @(8, 45917.0625, 0, 13.118, 0, 0, "17.09.20257"),
    @(9, 45917.07291666666, 0, 10.847, 0, 0, "17.09.20258"),
    @(10, 45917.08333333334, 0.05, 3.372, 0, 0, "17.09.20259"),
    @(11, 45917.09375, 0.009, 0.625, 0, 0, "17.09.202510"),
    @(12, 45917.10416666666, 0.116, 0.778, 0, 0, "17.09.202511"),
    @(13, 45917.11458333334, 0.073, 0.655, 0, 0, "17.09.202512"),
    @(14, 45917.125, 0.189, 0.28, 0, 0, "17.09.202513"),
    @(15, 45917.13541666666, 0.044, 1.669, 0, 0, "17.09.202514"),
    @(16, 45917.14583333334, 0.094, 1.323, 0, 0, "17.09.202515"),
    @(17, 45917.15625, 0.188, 0.139, 0, 0, "17.09.202516"),
    @(18, 45917.16666666666, 0.123, 0.202, 0, 0, "17.09.202517"),
    @(19, 45917.17708333334, 0.244, 0.069, 0, 0, "17.09.202518"),
    @(20, 45917.1875, 2.025, 0.065, 0, 0, "17.09.202519"),
    @(21, 45917.19791666666, 1.183, 0.059, 0, 0, "17.09.202520"),
    @(22, 45917.20833333334, 0, 3.864, 0, 0, "17.09.202521"),
    @(23, 45917.21875, 4.229, 0.049, 0, 0, "17.09.202522"),
    @(24, 45917.22916666666, 18.336, 0, 0, 0, "17.09.202523"),
    @(25, 45917.23958333334, 37.437, 0, 0, 0, "17.09.202524"),
    @(26, 45917.25, 17.098, 0, 0, 0, "17.09.202525"),
    @(27, 45917.26041666666, 40.798, 0, 0, 0, "17.09.202526"),
    @(28, 45917.27083333334, 51.023, 0, 37.5, 0, "17.09.202527"),
    @(29, 45917.28125, 27.205, 0, 62.5, 0, "17.09.202528"),
    @(30, 45917.29166666666, 35.705, 0, 62.5, 0, "17.09.202529"),
    @(31, 45917.30208333334, 68.56, 0, 62.5, 0, "17.09.202530"),
    @(32, 45917.3125, 56.074, 0, 112.5, 0, "17.09.202531"),
    @(33, 45917.32291666666, 4.334, 1.257, 150, 0, "17.09.202532"),
    @(34, 45917.33333333334, 25.868, 0.021, 37.5, 0, "17.09.202533"),
    @(35, 45917.34375, 36.52, 0, 37.5, 0, "17.09.202534"),
    @(36, 45917.35416666666, 15.659, 0, 50, 0, "17.09.202535"),
    @(37, 45917.36458333334, 5.385, 3.76, 50, 0, "17.09.202536"),
    @(38, 45917.375, 0, 16.222, 50, 0, "17.09.202537"),
    @(39, 45917.38541666666, 0, 13.877, 50, 0, "17.09.202538"),
    @(40, 45917.39583333334, 0, 20.608, 25, 0, "17.09.202539"),
    @(41, 45917.40625, 0, 56.224, 25, 0, "17.09.202540"),
    @(42, 45917.41666666666, 0, 20.609, 0, 0, "17.09.202541"),
    @(43, 45917.42708333334, 0, 33.701, 0, 0, "17.09.202542"),
    @(44, 45917.4375, 0, 50.599, 0, 0, "17.09.202543"),
    @(45, 45917.44791666666, 0, 49.975, 0, 0, "17.09.202544"),
    @(46, 45917.45833333334, 0.047, 6.912, 0, 0, "17.09.202545"),
    @(47, 45917.46875, 0.003, 6.991, 0, 0, "17.09.202546"),
    @(48, 45917.47916666666, 0, 14.468, 0, 0, "17.09.202547"),
    @(49, 45917.48958333334, 0, 29.737, 0, 0, "17.09.202548"),
    @(50, 45917.5, 0, 15.103, 0, 0, "17.09.202549"),
    @(51, 45917.51041666666, 0.014, 18.324, 0, 0, "17.09.202550"),
    @(52, 45917.52083333334, 0.309, 0.276, 0, 37.5, "17.09.202551"),
    @(53, 45917.53125, 0, 7.704, 0, 37.5, "17.09.202552"),
    @(54, 45917.54166666666, 0.142, 0.398, 0, 37.5, "17.09.202553"),
    @(55, 45917.55208333334, 0.093, 0.398, 0, 37.5, "17.09.202554"),
    @(56, 45917.5625, 0.16, 0.358, 0, 25, "17.09.202555"),
    @(57, 45917.57291666666, 0.006, 17.061, 0, 25, "17.09.202556"),
    @(58, 45917.58333333334, 0, 29.475, 0, 25, "17.09.202557"),
    @(59, 45917.59375, 0.006, 6.607, 0, 25, "17.09.202558"),
    @(60, 45917.60416666666, 3.446, 0.778, 0, 25, "17.09.202559"),
    @(61, 45917.61458333334, 3.676, 0.052, 0, 25, "17.09.202560"),
    @(62, 45917.625, 0, 21.869, 0, 0, "17.09.202561"),
    @(63, 45917.63541666666, 0.019, 4.314, 0, 25, "17.09.202562"),
    @(64, 45917.64583333334, 3.589, 0.032, 0, 25, "17.09.202563"),
    @(65, 45917.65625, 3.042, 0, 0, 25, "17.09.202564"),
    @(66, 45917.66666666666, 1.047, 3.904, 0, 0, "17.09.202565"),
    @(67, 45917.67708333334, 0.209, 0.933, 0, 0, "17.09.202566"),
    @(68, 45917.6875, 0.2, 0.207, 0, 0, "17.09.202567"),
    @(69, 45917.69791666666, 0.21, 0.232, 0, 0, "17.09.202568"),
    @(70, 45917.70833333334, 0, 28.217, 0, 0, "17.09.202569"),
    @(71, 45917.71875, 5.409, 0, 0, 0, "17.09.202570"),
    @(72, 45917.72916666666, 0.487, 0.044, 0, 0, "17.09.202571"),
    @(73, 45917.73958333334, 0.031, 0.839, 0, 0, "17.09.202572"),
    @(74, 45917.75, 0, 39.587, 0, 0, "17.09.202573"),
    @(75, 45917.76041666666, 0.024, 11.194, 0, 0, "17.09.202574"),
    @(76, 45917.77083333334, 0.318, 0.184, 0, 0, "17.09.202575"),
    @(77, 45917.78125, 1.743, 0, 0, 0, "17.09.202576"),
    @(78, 45917.79166666666, 12.794, 0, 0, 0, "17.09.202577"),
    @(79, 45917.80208333334, 36.167, 0, 0, 0, "17.09.202578"),
    @(80, 45917.8125, 41.921, 0, 0, 0, "17.09.202579"),
    @(81, 45917.82291666666, 41.323, 0, 0, 0, "17.09.202580"),
    @(82, 45917.83333333334, 45.623, 0, 0, 0, "17.09.202581"),
    @(83, 45917.84375, 5.082, 0.046, 48.25, 0, "17.09.202582"),
    @(84, 45917.85416666666, 0, 22.766, 48.25, 0, "17.09.202583"),
    @(85, 45917.86458333334, 0, 42.502, 0, 0, "17.09.202584"),
    @(86, 45917.875, 0.729, 8.823, 0, 0, "17.09.202585"),
    @(87, 45917.88541666666, 0, 16.946, 0, 0, "17.09.202586"),
    @(88, 45917.89583333334, 0, 29.265, 0, 0, "17.09.202587"),
    @(89, 45917.90625, 0, 53.604, 0, 0, "17.09.202588"),
    @(90, 45917.91666666666, 0.003, 15.944, 0, 0, "17.09.202589"),
    @(91, 45917.92708333334, 0, 16.642, 0, 0, "17.09.202590"),
    @(92, 45917.9375, 0, 21.47, 0, 0, "17.09.202591"),
    @(93, 45917.94791666666, 0, 27.019, 0, 0, "17.09.202592"),
    @(94, 45917.95833333334, 0.112, 3.426, 0, 0, "17.09.202593"),
    @(95, 45917.96875, 0.02, 0.276, 0, 0, "17.09.202594"),
    @(96, 45917.97916666666, 0, 15.447, 0, 0, "17.09.202595"),
    @(97, 45917.98958333334, 0, 50.938, 0, 0, "17.09.202596"),
    @(98, 45918, 0, 6.924, 0, 0, "18.09.20251"),
    @(99, 45918, 0, 6.924, 0, 0, "18.09.20251"),
    @(100, 45918.01041666666, 0, 6.386, 0, 0, "18.09.20252"),
    @(101, 45918.01041666666, 0, 6.386, 0, 0, "18.09.20252"),
    @(102, 45918.02083333334, 0, 20.396, 0, 0, "18.09.20253"),
    @(103, 45918.02083333334, 0, 20.396, 0, 0, "18.09.20253"),
    @(104, 45918.03125, 0, 26.501, 0, 0, "18.09.20254"),
    @(105, 45918.03125, 0, 26.501, 0, 0, "18.09.20254"),
    @(106, 45918.04166666666, 0, 24.222, 0, 0, "18.09.20255"),
    @(107, 45918.04166666666, 0, 24.222, 0, 0, "18.09.20255"),
    @(108, 45918.05208333334, 0, 19.788, 0, 0, "18.09.20256"),
    @(109, 45918.05208333334, 0, 19.788, 0, 0, "18.09.20256"),
    @(110, 45918.0625, 0, 36.503, 0, 0, "18.09.20257"),
    @(111, 45918.0625, 0, 36.503, 0, 0, "18.09.20257"),
    @(112, 45918.07291666666, 0, 30.799, 0, 0, "18.09.20258"),
    @(113, 45918.07291666666, 0, 30.799, 0, 0, "18.09.20258"),
    @(114, 45918.08333333334, 0, 27.876, 0, 0, "18.09.20259"),
    @(115, 45918.09375, 0, 18.105, 0, 0, "18.09.202510"),
    @(116, 45918.10416666666, 0, 5.994, 0, 0, "18.09.202511"),
    @(117, 45918.11458333334, 0, 17.033, 0, 0, "18.09.202512"),
    @(118, 45918.125, 0, 20.9, 0, 0, "18.09.202513"),
    @(119, 45918.13541666666, 0, 21.19, 0, 0, "18.09.202514"),
    @(120, 45918.14583333334, 0, 1.987, 0, 0, "18.09.202515"),
    @(121, 45918.15625, 0, 10.698, 0, 0, "18.09.202516"),
    @(122, 45918.16666666666, 0, 8.102, 0, 0, "18.09.202517"),
    @(123, 45918.17708333334, 0, 4.928, 0, 0, "18.09.202518"),
    @(124, 45918.1875, 0.948, 0.046, 0, 0, "18.09.202519"),
    @(125, 45918.19791666666, 0.949, 0.118, 0, 0, "18.09.202520"),
    @(126, 45918.20833333334, 0.047, 5.68, 0, 0, "18.09.202521"),
    @(127, 45918.21875, 0.001, 8.303, 0, 0, "18.09.202522"),
    @(128, 45918.22916666666, 0.386, 5.629, 0, 0, "18.09.202523"),
    @(129, 45918.23958333334, 0.115, 1.607, 0, 0, "18.09.202524"),
    @(130, 45918.25, 0, 21.02, 0, 0, "18.09.202525"),
    @(131, 45918.26041666666, 0.554, 1.932, 0, 0, "18.09.202526"),
    @(132, 45918.27083333334, 1.259, 0.133, 0, 0, "18.09.202527"),
    @(133, 45918.28125, 0.915, 0.138, 0, 0, "18.09.202528"),
    @(134, 45918.29166666666, 6.339, 0, 0, 0, "18.09.202529"),
    @(135, 45918.30208333334, 0.017, 15.522, 0, 0, "18.09.202530"),
    @(136, 45918.3125, 0, 50.923, 0, 0, "18.09.202531"),
    @(137, 45918.32291666666, 0, 85.378, 0, 0, "18.09.202532"),
    @(138, 45918.33333333334, 5.822, 8.662, 0, 0, "18.09.202533"),
    @(139, 45918.34375, 0, 36.821, 0, 0, "18.09.202534"),
    @(140, 45918.35416666666, 0, 78.594, 0, 0, "18.09.202535"),
    @(141, 45918.36458333334, 0, 78.918, 0, 0, "18.09.202536"),
    @(142, 45918.375, 0, 13.21, 0, 0, "18.09.202537"),
    @(143, 45918.38541666666, 0.016, 10.335, 0, 25, "18.09.202538"),
    @(144, 45918.39583333334, 0, 18.445, 0, 25, "18.09.202539"),
    @(145, 45918.40625, 0, 29.688, 0, 25, "18.09.202540"),
    @(146, 45918.41666666666, 2.249, 4.29, 0, 25, "18.09.202541"),
    @(147, 45918.42708333334, 0.18, 0.908, 0, 25, "18.09.202542")
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}

# New rows (146, 147) need F (Quarter) values explicitly since they did not exist before.
$ws.Cells.Item(146, 6).Value = 41
$ws.Cells.Item(147, 6).Value = 42
